$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4: was a blank space, now holds the text "100"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "100"
$ws.Range("A4").Style = "Normal"

# B4: was the numeric 60, now holds the blank-space text that A4 used to have
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = " "
$ws.Range("B4").Style = "Normal"

# C4: was a styled date serial number, now holds the literal text "2020-09-07"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2020-09-07"
$ws.Range("C4").Style = "Normal"
